$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "245.02"
Set-TextValue $ws.Range("G2") "17"

Set-TextValue $ws.Range("D3") "25.03"
Set-TextValue $ws.Range("G3") "17"

Set-TextValue $ws.Range("B4") "LEO"
Set-TextValue $ws.Range("C4") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D4") "3.501"
Set-TextValue $ws.Range("E4") "3LEOLEO"
Set-TextValue $ws.Range("G4") "17"

Set-TextValue $ws.Range("B5") "HuobiToken"
Set-TextValue $ws.Range("C5") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D5") "5.138"
Set-TextValue $ws.Range("E5") "4HuobiTokenHT"
Set-TextValue $ws.Range("G5") "17"

Set-TextValue $ws.Range("B6") "Cronos"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D6") "0.05634"
Set-TextValue $ws.Range("E6") "5CronosCRO"
Set-TextValue $ws.Range("G6") "17"

Set-TextValue $ws.Range("B7") "KuCoinToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D7") "6.530"
Set-TextValue $ws.Range("E7") "6KuCoinTokenKCS"
Set-TextValue $ws.Range("G7") "17"

Set-TextValue $ws.Range("B8") "GateToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D8") "2.981"
Set-TextValue $ws.Range("E8") "7GateTokenGT"
Set-TextValue $ws.Range("G8") "17"

Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.8115"
Set-TextValue $ws.Range("E9") "8MXTokenMX"
Set-TextValue $ws.Range("G9") "17"

Set-TextValue $ws.Range("B10") "FTXToken"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D10") "0.8376"
Set-TextValue $ws.Range("E10") "9FTXTokenFTT"
Set-TextValue $ws.Range("G10") "17"

Set-TextValue $ws.Range("B11") "One"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D11") "0.0005949"
Set-TextValue $ws.Range("E11") "10OneONE"
Set-TextValue $ws.Range("G11") "17"

Set-TextValue $ws.Range("B12") "WazirX"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D12") "0.1337"
Set-TextValue $ws.Range("E12") "11WazirXWRX"
Set-TextValue $ws.Range("G12") "17"

Set-TextValue $ws.Range("B13") "MandalaExchangeToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D13") "0.06954"
Set-TextValue $ws.Range("E13") "12MandalaExchangeTokenMDX"
Set-TextValue $ws.Range("G13") "17"

Set-TextValue $ws.Range("B14") "BitrueCoin"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.02845"
Set-TextValue $ws.Range("E14") "13BitrueCoinBTR"
Set-TextValue $ws.Range("G14") "17"

Set-TextValue $ws.Range("B15") "BitMartToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09419"
Set-TextValue $ws.Range("E15") "14BitMartTokenBMX"
Set-TextValue $ws.Range("G15") "17"

Set-TextValue $ws.Range("B16") "BitForexToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001525"
Set-TextValue $ws.Range("E16") "15BitForexTokenBF"
Set-TextValue $ws.Range("G16") "17"

Set-TextValue $ws.Range("B17") "TigerCash"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.006113"
Set-TextValue $ws.Range("E17") "16TigerCashTCH"
Set-TextValue $ws.Range("G17") "17"

Set-TextValue $ws.Range("G18") "17"

Set-TextValue $ws.Range("D19") "0.3170"
Set-TextValue $ws.Range("G19") "17"

Set-TextValue $ws.Range("D20") "0.03269"
Set-TextValue $ws.Range("G20") "17"

Set-TextValue $ws.Range("G21") "17"

Set-TextValue $ws.Range("D22") "3.740"
Set-TextValue $ws.Range("G22") "17"

Set-TextValue $ws.Range("D23") "0.04706"
Set-TextValue $ws.Range("G23") "17"

Set-TextValue $ws.Range("G24") "17"

Set-TextValue $ws.Range("D25") "0.001242"
Set-TextValue $ws.Range("G25") "17"

Set-TextValue $ws.Range("D26") "0.004527"
Set-TextValue $ws.Range("G26") "17"

Set-TextValue $ws.Range("E27") "26NitroExNTX"
Set-TextValue $ws.Range("G27") "17"

Set-TextValue $ws.Range("D28") "0.0001708"
Set-TextValue $ws.Range("G28") "17"

Set-TextValue $ws.Range("G29") "17"

Set-TextValue $ws.Range("G30") "17"

Set-TextValue $ws.Range("G31") "17"

Set-TextValue $ws.Range("G32") "17"

Set-TextValue $ws.Range("G33") "17"

Set-TextValue $ws.Range("G34") "17"

Set-TextValue $ws.Range("G35") "17"

Set-TextValue $ws.Range("G36") "17"

Set-TextValue $ws.Range("G37") "17"

Set-TextValue $ws.Range("G38") "17"

Set-TextValue $ws.Range("G39") "17"

Set-TextValue $ws.Range("D40") "0.03627"
Set-TextValue $ws.Range("G40") "17"

Set-TextValue $ws.Range("D41") "0.006230"
Set-TextValue $ws.Range("E41") "40KickTokenKICKBestin24h"
Set-TextValue $ws.Range("G41") "17"

Set-TextValue $ws.Range("D42") "0.1052"
Set-TextValue $ws.Range("G42") "17"

Set-TextValue $ws.Range("D43") "0.002468"
Set-TextValue $ws.Range("G43") "17"

Set-TextValue $ws.Range("D44") "0.008340"
Set-TextValue $ws.Range("G44") "17"

Set-TextValue $ws.Range("D45") "0.00005285"
Set-TextValue $ws.Range("G45") "17"

Set-TextValue $ws.Range("G46") "17"

Set-TextValue $ws.Range("D47") "0.2200"
Set-TextValue $ws.Range("G47") "17"

Set-TextValue $ws.Range("D48") "0.002287"
Set-TextValue $ws.Range("G48") "17"

Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("G49") "17"

Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("G50") "17"

Set-TextValue $ws.Range("G51") "17"
